# Helper: force a cell to hold a TEXT value (avoids Excel auto-converting
# numeric-looking strings like "011686" or "0.57" into numbers), while
# keeping the cell's style index at the default (no style) like the source
# data cells.
function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) "总计" (totals) sheet: insert a new row 2 for the 2022-Q3 summary
#    and bump the running index (column A) of every row below it.
# ------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

$total.Rows.Item(2).Insert()
$total.Range("B2:D2").ClearFormats()

# restore the bordered/bold style on the new A2 (matches the other
# index cells in column A) by copying the format from A3.
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.07

# shift the running index of the previously-existing rows down by one
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5
$total.Range("A8").Value = 6
$total.Range("A9").Value = 7

# ------------------------------------------------------------------
# 2) New "2022-Q3" detail sheet: clone the existing "2022-Q2" sheet
#    (same column layout/styles) right after "总计", rename it, trim
#    it down to 2 data rows, and fill in the 2022-Q3 numbers.
# ------------------------------------------------------------------
$template = $wb.Worksheets.Item(2)
$template.Copy($null, $total)

$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# the template had 4 data rows (rows 2-5); 2022-Q3 only needs 2 (rows 2-3)
$q3.Range("A4:H5").Delete(-4162)

# Row 2
$q3.Range("A2").Value = 0
Set-TextValue $q3.Range("B2") "011686"
$q3.Range("C2").Value = "创金合信先进装备股票C"
Set-TextValue $q3.Range("D2") "0.57"
Set-TextValue $q3.Range("E2") "80.17"
Set-TextValue $q3.Range("F2") "8.90"
Set-TextValue $q3.Range("G2") "0.0507"
$q3.Range("H2").Value = 5

# Row 3
$q3.Range("A3").Value = 1
Set-TextValue $q3.Range("B3") "011685"
$q3.Range("C3").Value = "创金合信先进装备股票A"
Set-TextValue $q3.Range("D3") "0.25"
Set-TextValue $q3.Range("E3") "80.17"
Set-TextValue $q3.Range("F3") "8.90"
Set-TextValue $q3.Range("G3") "0.0222"
$q3.Range("H3").Value = 5

Write-Output "done"
